# Update market/profit data values across sheets as part of a scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3812.8572
$ws.Range("I64").Value = 2995
$ws.Range("J64").Value = 4140
$ws.Range("K64").Value = 2995
$ws.Range("L64").Value = 4140
$ws.Range("M64").Value = -2747
$ws.Range("N64").Value = -4636
$ws.Range("H67").Value = 3812.8572
$ws.Range("I67").Value = 2995
$ws.Range("J67").Value = 4140
$ws.Range("K67").Value = 2995
$ws.Range("L67").Value = 4140
$ws.Range("M67").Value = -2137
$ws.Range("N67").Value = -5856
$ws.Range("H80").Value = 2151.9443
$ws.Range("I80").Value = 2067.1667
$ws.Range("J80").Value = 2194.3333
$ws.Range("K80").Value = 6201.500100000001
$ws.Range("L80").Value = 6582.999899999999
$ws.Range("M80").Value = -5203.500100000001
$ws.Range("N80").Value = -8578.999899999999
$ws.Range("H83").Value = 2151.9443
$ws.Range("I83").Value = 2067.1667
$ws.Range("J83").Value = 2194.3333
$ws.Range("K83").Value = 18604.5003
$ws.Range("L83").Value = 19748.9997
$ws.Range("M83").Value = -13612.5003
$ws.Range("N83").Value = -29732.9997
$ws.Range("H88").Value = 6067.353
$ws.Range("J88").Value = 7042.909
$ws.Range("L88").Value = 7042.909
$ws.Range("N88").Value = -7854.909
$ws.Range("H91").Value = 6067.353
$ws.Range("J91").Value = 7042.909
$ws.Range("L91").Value = 7042.909
$ws.Range("N91").Value = -9850.909
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("H106").Value = 1118271.5
$ws.Range("I106").Value = 3334547.8
$ws.Range("J106").Value = 10133.333
$ws.Range("K106").Value = 3334547.8
$ws.Range("L106").Value = 10133.333
$ws.Range("M106").Value = -3333916.8
$ws.Range("N106").Value = -11395.333
$ws.Range("H125").Value = 1166.6666
$ws.Range("J125").Value = 1166.6666
$ws.Range("L125").Value = 10499.9994
$ws.Range("N125").Value = -15419.9994
$ws.Range("H127").Value = 142863000
$ws.Range("I127").Value = 1000000000
$ws.Range("J127").Value = 6844.3335
$ws.Range("K127").Value = 3000000000
$ws.Range("L127").Value = 20533.0005
$ws.Range("M127").Value = -2999995040
$ws.Range("N127").Value = -30453.0005
$ws.Range("H130").Value = 20999.9
$ws.Range("J130").Value = 20999.9
$ws.Range("L130").Value = 20999.9
$ws.Range("N130").Value = -31039.9
$ws.Range("H132").Value = 4060.442
$ws.Range("I132").Value = 1758.7949
$ws.Range("J132").Value = 26501.5
$ws.Range("K132").Value = 5276.384700000001
$ws.Range("L132").Value = 79504.5
$ws.Range("M132").Value = -2746.384700000001
$ws.Range("N132").Value = -84564.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 18799.8
$ws.Range("J52").Value = 18799.8
$ws.Range("L52").Value = 18799.8
$ws.Range("N52").Value = -19435.8
$ws.Range("H137").Value = 41800
$ws.Range("J137").Value = 41800
$ws.Range("L137").Value = 41800
$ws.Range("N137").Value = -52000
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1791.8
$ws.Range("I86").Value = 1605.5883
$ws.Range("J86").Value = 2187.5
$ws.Range("K86").Value = 1605.5883
$ws.Range("L86").Value = 2187.5
$ws.Range("M86").Value = -482.5882999999999
$ws.Range("N86").Value = -4433.5
$ws.Range("H89").Value = 1791.8
$ws.Range("I89").Value = 1605.5883
$ws.Range("J89").Value = 2187.5
$ws.Range("K89").Value = 8027.941499999999
$ws.Range("L89").Value = 10937.5
$ws.Range("M89").Value = -2411.941499999999
$ws.Range("N89").Value = -22169.5
$ws.Range("H134").Value = 536109.5600000001
$ws.Range("I134").Value = 836021.5
$ws.Range("J134").Value = 2932.7036
$ws.Range("K134").Value = 2508064.5
$ws.Range("L134").Value = 8798.110799999999
$ws.Range("M134").Value = -2505529.5
$ws.Range("N134").Value = -13868.1108
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1289.975
$ws.Range("I31").Value = 837.5599999999999
$ws.Range("J31").Value = 2044
$ws.Range("K31").Value = 837.5599999999999
$ws.Range("L31").Value = 2044
$ws.Range("M31").Value = -542.5599999999999
$ws.Range("N31").Value = -2634
$ws.Range("H34").Value = 1289.975
$ws.Range("I34").Value = 837.5599999999999
$ws.Range("J34").Value = 2044
$ws.Range("K34").Value = 837.5599999999999
$ws.Range("L34").Value = 2044
$ws.Range("M34").Value = -635.5599999999999
$ws.Range("N34").Value = -2448
$ws.Range("H58").Value = 4547.7334
$ws.Range("I58").Value = 5009.769
$ws.Range("J58").Value = 1544.5
$ws.Range("K58").Value = 5009.769
$ws.Range("L58").Value = 1544.5
$ws.Range("M58").Value = -4806.769
$ws.Range("N58").Value = -1950.5
$ws.Range("H132").Value = 1236938.4
$ws.Range("I132").Value = 2482.5908
$ws.Range("K132").Value = 7447.7724
$ws.Range("M132").Value = -4917.7724
$ws.Range("H136").Value = 4547.7334
$ws.Range("I136").Value = 5009.769
$ws.Range("J136").Value = 1544.5
$ws.Range("K136").Value = 15029.307
$ws.Range("L136").Value = 4633.5
$ws.Range("M136").Value = -12479.307
$ws.Range("N136").Value = -9733.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 663
$ws.Range("J80").Value = 663
$ws.Range("L80").Value = 1989
$ws.Range("N80").Value = -3861
$ws.Range("H83").Value = 663
$ws.Range("J83").Value = 663
$ws.Range("L83").Value = 5967
$ws.Range("N83").Value = -15327
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 22687.375
$ws.Range("J63").Value = 22857.143
$ws.Range("L63").Value = 22857.143
$ws.Range("N63").Value = -24229.143
$ws.Range("H66").Value = 22687.375
$ws.Range("J66").Value = 22857.143
$ws.Range("L66").Value = 68571.429
$ws.Range("N66").Value = -75435.429
$ws.Range("H70").Value = 5223.3335
$ws.Range("I70").Value = 4885.3335
$ws.Range("J70").Value = 5786.6665
$ws.Range("K70").Value = 4885.3335
$ws.Range("L70").Value = 5786.6665
$ws.Range("M70").Value = -4615.3335
$ws.Range("N70").Value = -6326.6665
$ws.Range("H73").Value = 5223.3335
$ws.Range("I73").Value = 4885.3335
$ws.Range("J73").Value = 5786.6665
$ws.Range("K73").Value = 4885.3335
$ws.Range("L73").Value = 5786.6665
$ws.Range("M73").Value = -3949.3335
$ws.Range("N73").Value = -7658.6665
$ws.Range("H80").Value = 3150
$ws.Range("I80").Value = 2780
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 2780
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -1782
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 3150
$ws.Range("I83").Value = 2780
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 13900
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -8908
$ws.Range("N83").Value = -34984
$ws.Range("H132").Value = 2224945.5
$ws.Range("I132").Value = 2923.8572
$ws.Range("J132").Value = 5884745.5
$ws.Range("K132").Value = 8771.571599999999
$ws.Range("L132").Value = 17654236.5
$ws.Range("M132").Value = -6241.571599999999
$ws.Range("N132").Value = -17659296.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 33500
$ws.Range("J64").Value = 33500
$ws.Range("L64").Value = 33500
$ws.Range("N64").Value = -33950
$ws.Range("H67").Value = 33500
$ws.Range("J67").Value = 33500
$ws.Range("L67").Value = 33500
$ws.Range("N67").Value = -35060
$ws.Range("H132").Value = 2963
$ws.Range("I132").Value = 2960.7805
$ws.Range("J132").Value = 2973.111
$ws.Range("K132").Value = 8882.341499999999
$ws.Range("L132").Value = 8919.332999999999
$ws.Range("M132").Value = -6352.341499999999
$ws.Range("N132").Value = -13979.333
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 485.6
$ws.Range("I113").Value = 437.23077
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1311.69231
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 858.3076900000001
$ws.Range("N113").Value = -6740
$ws.Range("H136").Value = 1334.5
$ws.Range("I136").Value = 753.30554
$ws.Range("J136").Value = 2206.2917
$ws.Range("K136").Value = 2259.91662
$ws.Range("L136").Value = 6618.875100000001
$ws.Range("M136").Value = 290.08338
$ws.Range("N136").Value = -11718.8751
